$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on column C (the "Förändrad" / Changed date column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# Update column C (rows 2..lastRow) from 46061 to 46062 (one day later)
$ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3)).Value = 46062
